$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BACKLOG")

# Update "data do BD" (K column) for several backlog items, which flips
# their "concluído" flag (J column formula) to TRUE and recalculates the
# burndown summary (O2:Q7) plus the chart cache automatically.
$ws.Range("K23").Value = 45963
$ws.Range("K26").Value = 45963
$ws.Range("K33").Value = 45962
$ws.Range("K34").Value = 45963
$ws.Range("K35").Value = 45963
$ws.Range("K43").Value = 45995
$ws.Range("K44").Value = 45995

# Reclassify row 24's priority from IMPORTANTE to DESEJÁVEL.
$ws.Range("E24").Value = "DESEJÁVEL"

# Update the visible/selected cell & scroll position to match the saved view.
$ws.Range("K24").Select()

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
